# Insert a new data row before row 38 (shifts old rows 38-153 down to 39-154)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("38").Insert()

# Populate the new row 38 with its data (same market/category as the rest of the sheet)
$ws.Range("A38").Value = 10
$ws.Range("B38").Value = "Vega Modelo de Temuco"
$ws.Range("C38").Value = "La Araucanía"
$ws.Range("D38").Value = 44526
$ws.Range("E38").Value = 9
$ws.Range("F38").Value = 100112052
$ws.Range("G38").Value = "Albahaca"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 20
$ws.Range("K38").Value = 3500
$ws.Range("L38").Value = 3500
$ws.Range("M38").Value = 3500
$ws.Range("N38").Value = "$/paquete"
$ws.Range("O38").Value = "Región del Maule"
$ws.Range("P38").Value = 3500
$ws.Range("Q38").Value = 1
$ws.Range("R38").Value = "Hortaliza"
